$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - "It would succeed in limiting climate change"
$ws.Range("A2").Value = "It would succeed in limiting climate change"
$ws.Range("B2").Value = 0.672053541001286
$ws.Range("C2").Value = 0.582938067961184
$ws.Range("D2").Value = 0.540830920916378
$ws.Range("E2").Value = 0.656789582949536
$ws.Range("F2").Value = 0.538104795906038
$ws.Range("G2").Value = 0.541350478300235

# Row 3 - "It would hurt the [Country] economy"
$ws.Range("A3").Value = "It would hurt the [Country] economy"
$ws.Range("B3").Value = 0.891301757279309
$ws.Range("C3").Value = 0.780204464912827
$ws.Range("D3").Value = 0.659795214245632
$ws.Range("E3").Value = 0.729457954248223
$ws.Range("F3").Value = 0.73782451273014
$ws.Range("G3").Value = 0.902695682752281

# Row 4 - "It would penalize my household"
$ws.Range("A4").Value = "It would penalize my household"
$ws.Range("B4").Value = 0.839128610889381
$ws.Range("C4").Value = 0.742901276441401
$ws.Range("D4").Value = 0.7017279146792
$ws.Range("E4").Value = 0.728815269477337
$ws.Range("F4").Value = 0.750395530604615
$ws.Range("G4").Value = 0.782192072002497

# Row 5 - "It would make people change their lifestyle"
$ws.Range("A5").Value = "It would make people change their lifestyle"
$ws.Range("B5").Value = 0.724595557495586
$ws.Range("C5").Value = 0.626578181406864
$ws.Range("D5").Value = 0.658937842683342
$ws.Range("E5").Value = 0.644915946352901
$ws.Range("F5").Value = 0.597782105482217
$ws.Range("G5").Value = 0.580143827052933

# Row 6 - "It would reduce poverty in low-income countries"
$ws.Range("A6").Value = "It would reduce poverty in low-income countries"
$ws.Range("B6").Value = 0.689054219802248
$ws.Range("C6").Value = 0.664380381636305
$ws.Range("D6").Value = 0.704094709238634
$ws.Range("E6").Value = 0.711275867180752
$ws.Range("F6").Value = 0.572528055051728
$ws.Range("G6").Value = 0.60927373973101

# Row 7 - "It might be detrimental to some poor countries"
$ws.Range("A7").Value = "It might be detrimental to some poor countries"
$ws.Range("B7").Value = 0.74503010611439
$ws.Range("C7").Value = 0.652681493456455
$ws.Range("D7").Value = 0.646895029917251
$ws.Range("E7").Value = 0.659443297789059
$ws.Range("F7").Value = 0.708152559274898
$ws.Range("G7").Value = 0.636222901700506

# Row 8 - "It could foster global cooperation"
$ws.Range("A8").Value = "It could foster global cooperation"
$ws.Range("B8").Value = 0.719660543344661
$ws.Range("C8").Value = 0.573814842079329
$ws.Range("D8").Value = 0.523468296219712
$ws.Range("E8").Value = 0.634321205004898
$ws.Range("F8").Value = 0.525018018015189
$ws.Range("G8").Value = 0.562108110779417

# Row 9 - "It could fuel corruption in low-income countries"
$ws.Range("A9").Value = "It could fuel corruption in low-income countries"
$ws.Range("B9").Value = 0.77708533855146
$ws.Range("C9").Value = 0.706642015509628
$ws.Range("D9").Value = 0.782414898061268
$ws.Range("E9").Value = 0.655380261897868
$ws.Range("F9").Value = 0.777203628123294
$ws.Range("G9").Value = 0.689685291309682

# Row 10 - "It could be subject to fraud"
$ws.Range("A10").Value = "It could be subject to fraud"
$ws.Range("B10").Value = 0.834418003882785
$ws.Range("C10").Value = 0.798367178476084
$ws.Range("D10").Value = 0.775382352209707
$ws.Range("E10").Value = 0.796236245924259
$ws.Range("F10").Value = 0.817054675660094
$ws.Range("G10").Value = 0.794298090146449

# Row 11 - "It would be technically difficult to put in place"
$ws.Range("A11").Value = "It would be technically difficult to put in place"
$ws.Range("B11").Value = 0.811160240941621
$ws.Range("C11").Value = 0.71958710509228
$ws.Range("D11").Value = 0.697265999422906
$ws.Range("E11").Value = 0.668537258017323
$ws.Range("F11").Value = 0.73485497102585
$ws.Range("G11").Value = 0.740334261560751

# Row 12 - "Having enough information on\nthis scheme and its consequences"
$ws.Range("A12").Value = "Having enough information on`nthis scheme and its consequences"
$ws.Range("B12").Value = 0.909620188686539
$ws.Range("C12").Value = 0.765524000555159
$ws.Range("D12").Value = 0.784696947296919
$ws.Range("E12").Value = 0.712068192462401
$ws.Range("F12").Value = 0.792464056175893
$ws.Range("G12").Value = 0.791230367967449
